# Auto-generated edit script: refresh scraped schedule data for
# sheets LP1912 (sheet1), LP1912-215 (sheet2) and 6203-6173 (sheet3).
# Mirrors a new scrape run at 19:11:59 that appended newly observed
# arrivals (sorted ascending by Hora_Llegada) and bumped the header
# timestamps / row counts on all three sheets.

$wb = $excel.ActiveWorkbook

$newScrapeTime = "19:11:59"
$updatedLabel = "Última actualización: " + $newScrapeTime

# ---------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = $updatedLabel
$ws1.Range("A3").Value = "Total filas: 111"

$data1 = New-Object 'object[,]' 111,5
$data1[0,0] = "16:46:42"
$data1[0,1] = "16:47"
$data1[0,2] = "15_ABASTO"
$data1[0,3] = 1
$data1[0,4] = "LP1912"
$data1[1,0] = "16:50:41"
$data1[1,1] = "16:50"
$data1[1,2] = "10_OLMOS"
$data1[1,3] = 0
$data1[1,4] = "LP1912"
$data1[2,0] = "16:50:41"
$data1[2,1] = "16:51"
$data1[2,2] = "15_ABASTO"
$data1[2,3] = 1
$data1[2,4] = "LP1912"
$data1[3,0] = "16:52:37"
$data1[3,1] = "16:53"
$data1[3,2] = "15_ABASTO"
$data1[3,3] = 1
$data1[3,4] = "LP1912"
$data1[4,0] = "16:46:42"
$data1[4,1] = "16:53"
$data1[4,2] = "10_OLMOS"
$data1[4,3] = 7
$data1[4,4] = "LP1912"
$data1[5,0] = "16:46:42"
$data1[5,1] = "16:56"
$data1[5,2] = "215C_EL PATO"
$data1[5,3] = 10
$data1[5,4] = "LP1912"
$data1[6,0] = "16:46:42"
$data1[6,1] = "17:01"
$data1[6,2] = "16_SANTA ANA"
$data1[6,3] = 15
$data1[6,4] = "LP1912"
$data1[7,0] = "16:46:42"
$data1[7,1] = "17:03"
$data1[7,2] = "23_HERNANDEZ"
$data1[7,3] = 17
$data1[7,4] = "LP1912"
$data1[8,0] = "16:46:42"
$data1[8,1] = "17:04"
$data1[8,2] = "14_ABASTO"
$data1[8,3] = 18
$data1[8,4] = "LP1912"
$data1[9,0] = "16:46:42"
$data1[9,1] = "17:07"
$data1[9,2] = "15_ABASTO"
$data1[9,3] = 21
$data1[9,4] = "LP1912"
$data1[10,0] = "16:46:42"
$data1[10,1] = "17:13"
$data1[10,2] = "23_HERNANDEZ"
$data1[10,3] = 27
$data1[10,4] = "LP1912"
$data1[11,0] = "16:46:42"
$data1[11,1] = "17:14"
$data1[11,2] = "10_OLMOS"
$data1[11,3] = 28
$data1[11,4] = "LP1912"
$data1[12,0] = "17:13:30"
$data1[12,1] = "17:16"
$data1[12,2] = "10_OLMOS"
$data1[12,3] = 3
$data1[12,4] = "LP1912"
$data1[13,0] = "16:46:42"
$data1[13,1] = "17:17"
$data1[13,2] = "17_ROMERO"
$data1[13,3] = 31
$data1[13,4] = "LP1912"
$data1[14,0] = "16:50:41"
$data1[14,1] = "17:17"
$data1[14,2] = "23_HERNANDEZ"
$data1[14,3] = 27
$data1[14,4] = "LP1912"
$data1[15,0] = "16:52:37"
$data1[15,1] = "17:20"
$data1[15,2] = "23_HERNANDEZ"
$data1[15,3] = 28
$data1[15,4] = "LP1912"
$data1[16,0] = "16:46:42"
$data1[16,1] = "17:23"
$data1[16,2] = "16_SANTA ANA"
$data1[16,3] = 37
$data1[16,4] = "LP1912"
$data1[17,0] = "16:46:42"
$data1[17,1] = "17:24"
$data1[17,2] = "11_ETCHEVERRY"
$data1[17,3] = 38
$data1[17,4] = "LP1912"
$data1[18,0] = "17:13:30"
$data1[18,1] = "17:27"
$data1[18,2] = "15_ABASTO"
$data1[18,3] = 14
$data1[18,4] = "LP1912"
$data1[19,0] = "17:13:30"
$data1[19,1] = "17:33"
$data1[19,2] = "23_HERNANDEZ"
$data1[19,3] = 20
$data1[19,4] = "LP1912"
$data1[20,0] = "16:50:41"
$data1[20,1] = "17:34"
$data1[20,2] = "10_OLMOS"
$data1[20,3] = 44
$data1[20,4] = "LP1912"
$data1[21,0] = "16:46:42"
$data1[21,1] = "17:35"
$data1[21,2] = "16_P MOR-SANTA ANA"
$data1[21,3] = 49
$data1[21,4] = "LP1912"
$data1[22,0] = "16:52:37"
$data1[22,1] = "17:36"
$data1[22,2] = "27_EL RETIRO"
$data1[22,3] = 44
$data1[22,4] = "LP1912"
$data1[23,0] = "17:13:30"
$data1[23,1] = "17:37"
$data1[23,2] = "27_EL RETIRO"
$data1[23,3] = 24
$data1[23,4] = "LP1912"
$data1[24,0] = "16:46:42"
$data1[24,1] = "17:38"
$data1[24,2] = "17X38_ROMERO"
$data1[24,3] = 52
$data1[24,4] = "LP1912"
$data1[25,0] = "17:35:09"
$data1[25,1] = "17:39"
$data1[25,2] = "27_EL RETIRO"
$data1[25,3] = 4
$data1[25,4] = "LP1912"
$data1[26,0] = "17:13:30"
$data1[26,1] = "17:41"
$data1[26,2] = "23_HERNANDEZ"
$data1[26,3] = 28
$data1[26,4] = "LP1912"
$data1[27,0] = "16:46:42"
$data1[27,1] = "17:44"
$data1[27,2] = "215B_EL PATO"
$data1[27,3] = 58
$data1[27,4] = "LP1912"
$data1[28,0] = "17:35:09"
$data1[28,1] = "17:45"
$data1[28,2] = "215B_EL PATO"
$data1[28,3] = 10
$data1[28,4] = "LP1912"
$data1[29,0] = "16:50:41"
$data1[29,1] = "17:47"
$data1[29,2] = "16_SANTA ANA"
$data1[29,3] = 57
$data1[29,4] = "LP1912"
$data1[30,0] = "16:46:42"
$data1[30,1] = "17:48"
$data1[30,2] = "27_EL RETIRO"
$data1[30,3] = 62
$data1[30,4] = "LP1912"
$data1[31,0] = "16:50:41"
$data1[31,1] = "17:49"
$data1[31,2] = "27_EL RETIRO"
$data1[31,3] = 59
$data1[31,4] = "LP1912"
$data1[32,0] = "16:46:42"
$data1[32,1] = "17:50"
$data1[32,2] = "215_EL PELIGRO"
$data1[32,3] = 64
$data1[32,4] = "LP1912"
$data1[33,0] = "17:47:22"
$data1[33,1] = "17:51"
$data1[33,2] = "215B_EL PATO"
$data1[33,3] = 4
$data1[33,4] = "LP1912"
$data1[34,0] = "16:52:37"
$data1[34,1] = "17:51"
$data1[34,2] = "215_EL PELIGRO"
$data1[34,3] = 59
$data1[34,4] = "LP1912"
$data1[35,0] = "17:47:22"
$data1[35,1] = "17:54"
$data1[35,2] = "10_OLMOS"
$data1[35,3] = 7
$data1[35,4] = "LP1912"
$data1[36,0] = "17:35:09"
$data1[36,1] = "18:00"
$data1[36,2] = "16_SANTA ANA"
$data1[36,3] = 25
$data1[36,4] = "LP1912"
$data1[37,0] = "16:46:42"
$data1[37,1] = "18:02"
$data1[37,2] = "17_ROMERO"
$data1[37,3] = 76
$data1[37,4] = "LP1912"
$data1[38,0] = "16:52:37"
$data1[38,1] = "18:03"
$data1[38,2] = "17_ROMERO"
$data1[38,3] = 71
$data1[38,4] = "LP1912"
$data1[39,0] = "17:35:09"
$data1[39,1] = "18:03"
$data1[39,2] = "23_HERNANDEZ"
$data1[39,3] = 28
$data1[39,4] = "LP1912"
$data1[40,0] = "16:46:42"
$data1[40,1] = "18:04"
$data1[40,2] = "14_ABASTO"
$data1[40,3] = 78
$data1[40,4] = "LP1912"
$data1[41,0] = "17:35:09"
$data1[41,1] = "18:05"
$data1[41,2] = "14_ABASTO"
$data1[41,3] = 30
$data1[41,4] = "LP1912"
$data1[42,0] = "18:10:41"
$data1[42,1] = "18:11"
$data1[42,2] = "16_SANTA ANA"
$data1[42,3] = 1
$data1[42,4] = "LP1912"
$data1[43,0] = "18:10:41"
$data1[43,1] = "18:11"
$data1[43,2] = "10_OLMOS"
$data1[43,3] = 1
$data1[43,4] = "LP1912"
$data1[44,0] = "16:52:37"
$data1[44,1] = "18:14"
$data1[44,2] = "10_OLMOS"
$data1[44,3] = 82
$data1[44,4] = "LP1912"
$data1[45,0] = "17:47:22"
$data1[45,1] = "18:21"
$data1[45,2] = "16_SANTA ANA"
$data1[45,3] = 34
$data1[45,4] = "LP1912"
$data1[46,0] = "16:46:42"
$data1[46,1] = "18:24"
$data1[46,2] = "11_ETCHEVERRY"
$data1[46,3] = 98
$data1[46,4] = "LP1912"
$data1[47,0] = "17:35:09"
$data1[47,1] = "18:25"
$data1[47,2] = "11_ETCHEVERRY"
$data1[47,3] = 50
$data1[47,4] = "LP1912"
$data1[48,0] = "17:13:30"
$data1[48,1] = "18:27"
$data1[48,2] = "15_ABASTO"
$data1[48,3] = 74
$data1[48,4] = "LP1912"
$data1[49,0] = "17:47:22"
$data1[49,1] = "18:31"
$data1[49,2] = "16_SANTA ANA"
$data1[49,3] = 44
$data1[49,4] = "LP1912"
$data1[50,0] = "17:35:09"
$data1[50,1] = "18:31"
$data1[50,2] = "23_HERNANDEZ"
$data1[50,3] = 56
$data1[50,4] = "LP1912"
$data1[51,0] = "17:54:43"
$data1[51,1] = "18:33"
$data1[51,2] = "23_HERNANDEZ"
$data1[51,3] = 39
$data1[51,4] = "LP1912"
$data1[52,0] = "16:46:42"
$data1[52,1] = "18:34"
$data1[52,2] = "14X44_ABASTO"
$data1[52,3] = 108
$data1[52,4] = "LP1912"
$data1[53,0] = "16:46:42"
$data1[53,1] = "18:38"
$data1[53,2] = "17X38_ROMERO"
$data1[53,3] = 112
$data1[53,4] = "LP1912"
$data1[54,0] = "16:46:42"
$data1[54,1] = "18:41"
$data1[54,2] = "16_P MOR-SANTA ANA"
$data1[54,3] = 115
$data1[54,4] = "LP1912"
$data1[55,0] = "17:13:30"
$data1[55,1] = "18:41"
$data1[55,2] = "14_ABASTO"
$data1[55,3] = 88
$data1[55,4] = "LP1912"
$data1[56,0] = "17:47:22"
$data1[56,1] = "18:44"
$data1[56,2] = "14_ABASTO"
$data1[56,3] = 57
$data1[56,4] = "LP1912"
$data1[57,0] = "17:35:09"
$data1[57,1] = "18:45"
$data1[57,2] = "14_ABASTO"
$data1[57,3] = 70
$data1[57,4] = "LP1912"
$data1[58,0] = "18:44:34"
$data1[58,1] = "18:47"
$data1[58,2] = "14_ABASTO"
$data1[58,3] = 3
$data1[58,4] = "LP1912"
$data1[59,0] = "17:35:09"
$data1[59,1] = "18:51"
$data1[59,2] = "15_ABASTO"
$data1[59,3] = 76
$data1[59,4] = "LP1912"
$data1[60,0] = "18:52:04"
$data1[60,1] = "18:52"
$data1[60,2] = "15_ABASTO"
$data1[60,3] = 0
$data1[60,4] = "LP1912"
$data1[61,0] = "17:54:43"
$data1[61,1] = "18:53"
$data1[61,2] = "16_SANTA ANA"
$data1[61,3] = 59
$data1[61,4] = "LP1912"
$data1[62,0] = "17:35:09"
$data1[62,1] = "18:59"
$data1[62,2] = "10_OLMOS"
$data1[62,3] = 84
$data1[62,4] = "LP1912"
$data1[63,0] = "17:13:30"
$data1[63,1] = "19:01"
$data1[63,2] = "17_ROMERO"
$data1[63,3] = 108
$data1[63,4] = "LP1912"
$data1[64,0] = "18:10:41"
$data1[64,1] = "19:03"
$data1[64,2] = "23_HERNANDEZ"
$data1[64,3] = 53
$data1[64,4] = "LP1912"
$data1[65,0] = "18:52:04"
$data1[65,1] = "19:04"
$data1[65,2] = "23_HERNANDEZ"
$data1[65,3] = 12
$data1[65,4] = "LP1912"
$data1[66,0] = "18:31:18"
$data1[66,1] = "19:05"
$data1[66,2] = "16_SANTA ANA"
$data1[66,3] = 34
$data1[66,4] = "LP1912"
$data1[67,0] = "17:13:30"
$data1[67,1] = "19:11"
$data1[67,2] = "81_EL PELIGRO"
$data1[67,3] = 118
$data1[67,4] = "LP1912"
$data1[68,0] = "18:10:41"
$data1[68,1] = "19:14"
$data1[68,2] = "14_ABASTO"
$data1[68,3] = 64
$data1[68,4] = "LP1912"
$data1[69,0] = "18:52:04"
$data1[69,1] = "19:15"
$data1[69,2] = "14_ABASTO"
$data1[69,3] = 23
$data1[69,4] = "LP1912"
$data1[70,0] = "17:47:22"
$data1[70,1] = "19:17"
$data1[70,2] = "27_EL RETIRO"
$data1[70,3] = 90
$data1[70,4] = "LP1912"
$data1[71,0] = "18:44:34"
$data1[71,1] = "19:17"
$data1[71,2] = "16_SANTA ANA"
$data1[71,3] = 33
$data1[71,4] = "LP1912"
$data1[72,0] = "17:35:09"
$data1[72,1] = "19:19"
$data1[72,2] = "27_EL RETIRO"
$data1[72,3] = 104
$data1[72,4] = "LP1912"
$data1[73,0] = "17:54:43"
$data1[73,1] = "19:20"
$data1[73,2] = "215C_EL PATO"
$data1[73,3] = 86
$data1[73,4] = "LP1912"
$data1[74,0] = "17:35:09"
$data1[74,1] = "19:21"
$data1[74,2] = "215C_EL PATO"
$data1[74,3] = 106
$data1[74,4] = "LP1912"
$data1[75,0] = "19:11:59"
$data1[75,1] = "19:22"
$data1[75,2] = "27_EL RETIRO"
$data1[75,3] = 11
$data1[75,4] = "LP1912"
$data1[76,0] = "18:44:34"
$data1[76,1] = "19:26"
$data1[76,2] = "27_EL RETIRO"
$data1[76,3] = 42
$data1[76,4] = "LP1912"
$data1[77,0] = "19:11:59"
$data1[77,1] = "19:28"
$data1[77,2] = "225_GOMEZ"
$data1[77,3] = 17
$data1[77,4] = "LP1912"
$data1[78,0] = "17:35:09"
$data1[78,1] = "19:29"
$data1[78,2] = "225_GOMEZ"
$data1[78,3] = 114
$data1[78,4] = "LP1912"
$data1[79,0] = "17:54:43"
$data1[79,1] = "19:30"
$data1[79,2] = "27_EL RETIRO"
$data1[79,3] = 96
$data1[79,4] = "LP1912"
$data1[80,0] = "17:54:43"
$data1[80,1] = "19:30"
$data1[80,2] = "215_EL PELIGRO"
$data1[80,3] = 96
$data1[80,4] = "LP1912"
$data1[81,0] = "17:35:09"
$data1[81,1] = "19:31"
$data1[81,2] = "215_EL PELIGRO"
$data1[81,3] = 116
$data1[81,4] = "LP1912"
$data1[82,0] = "18:10:41"
$data1[82,1] = "19:31"
$data1[82,2] = "27_EL RETIRO"
$data1[82,3] = 81
$data1[82,4] = "LP1912"
$data1[83,0] = "18:44:34"
$data1[83,1] = "19:33"
$data1[83,2] = "23_HERNANDEZ"
$data1[83,3] = 49
$data1[83,4] = "LP1912"
$data1[84,0] = "18:31:18"
$data1[84,1] = "19:34"
$data1[84,2] = "23_HERNANDEZ"
$data1[84,3] = 63
$data1[84,4] = "LP1912"
$data1[85,0] = "19:11:59"
$data1[85,1] = "19:38"
$data1[85,2] = "16_SANTA ANA"
$data1[85,3] = 27
$data1[85,4] = "LP1912"
$data1[86,0] = "19:11:59"
$data1[86,1] = "19:39"
$data1[86,2] = "17X38_ROMERO"
$data1[86,3] = 28
$data1[86,4] = "LP1912"
$data1[87,0] = "19:11:59"
$data1[87,1] = "19:40"
$data1[87,2] = "16_SANTA ANA"
$data1[87,3] = 29
$data1[87,4] = "LP1912"
$data1[88,0] = "17:47:22"
$data1[88,1] = "19:40"
$data1[88,2] = "17X38_ROMERO"
$data1[88,3] = 113
$data1[88,4] = "LP1912"
$data1[89,0] = "19:11:59"
$data1[89,1] = "19:43"
$data1[89,2] = "11_ETCHEVERRY"
$data1[89,3] = 32
$data1[89,4] = "LP1912"
$data1[90,0] = "17:47:22"
$data1[90,1] = "19:44"
$data1[90,2] = "11_ETCHEVERRY"
$data1[90,3] = 117
$data1[90,4] = "LP1912"
$data1[91,0] = "18:31:18"
$data1[91,1] = "19:46"
$data1[91,2] = "11_ETCHEVERRY"
$data1[91,3] = 75
$data1[91,4] = "LP1912"
$data1[92,0] = "19:11:59"
$data1[92,1] = "19:50"
$data1[92,2] = "81_EL PELIGRO"
$data1[92,3] = 39
$data1[92,4] = "LP1912"
$data1[93,0] = "17:54:43"
$data1[93,1] = "19:51"
$data1[93,2] = "81_EL PELIGRO"
$data1[93,3] = 117
$data1[93,4] = "LP1912"
$data1[94,0] = "18:10:41"
$data1[94,1] = "19:58"
$data1[94,2] = "14X44_ABASTO"
$data1[94,3] = 108
$data1[94,4] = "LP1912"
$data1[95,0] = "18:31:18"
$data1[95,1] = "19:59"
$data1[95,2] = "14X44_ABASTO"
$data1[95,3] = 88
$data1[95,4] = "LP1912"
$data1[96,0] = "18:10:41"
$data1[96,1] = "20:00"
$data1[96,2] = "215C_EL PATO"
$data1[96,3] = 110
$data1[96,4] = "LP1912"
$data1[97,0] = "18:31:18"
$data1[97,1] = "20:01"
$data1[97,2] = "215C_EL PATO"
$data1[97,3] = 90
$data1[97,4] = "LP1912"
$data1[98,0] = "19:11:59"
$data1[98,1] = "20:04"
$data1[98,2] = "23_HERNANDEZ"
$data1[98,3] = 53
$data1[98,4] = "LP1912"
$data1[99,0] = "19:11:59"
$data1[99,1] = "20:13"
$data1[99,2] = "11_ETCHEVERRY"
$data1[99,3] = 62
$data1[99,4] = "LP1912"
$data1[100,0] = "18:31:18"
$data1[100,1] = "20:14"
$data1[100,2] = "11_ETCHEVERRY"
$data1[100,3] = 103
$data1[100,4] = "LP1912"
$data1[101,0] = "19:11:59"
$data1[101,1] = "20:25"
$data1[101,2] = "15_ABASTO"
$data1[101,3] = 74
$data1[101,4] = "LP1912"
$data1[102,0] = "18:31:18"
$data1[102,1] = "20:26"
$data1[102,2] = "15_ABASTO"
$data1[102,3] = 115
$data1[102,4] = "LP1912"
$data1[103,0] = "18:44:34"
$data1[103,1] = "20:28"
$data1[103,2] = "10_OLMOS"
$data1[103,3] = 104
$data1[103,4] = "LP1912"
$data1[104,0] = "18:31:18"
$data1[104,1] = "20:29"
$data1[104,2] = "10_OLMOS"
$data1[104,3] = 118
$data1[104,4] = "LP1912"
$data1[105,0] = "19:11:59"
$data1[105,1] = "20:43"
$data1[105,2] = "215B_EL PATO"
$data1[105,3] = 92
$data1[105,4] = "LP1912"
$data1[106,0] = "19:11:59"
$data1[106,1] = "20:44"
$data1[106,2] = "17X38_ROMERO"
$data1[106,3] = 93
$data1[106,4] = "LP1912"
$data1[107,0] = "18:52:04"
$data1[107,1] = "20:44"
$data1[107,2] = "215B_EL PATO"
$data1[107,3] = 112
$data1[107,4] = "LP1912"
$data1[108,0] = "18:52:04"
$data1[108,1] = "20:45"
$data1[108,2] = "17X38_ROMERO"
$data1[108,3] = 113
$data1[108,4] = "LP1912"
$data1[109,0] = "19:11:59"
$data1[109,1] = "21:01"
$data1[109,2] = "215A_EL PATO"
$data1[109,3] = 110
$data1[109,4] = "LP1912"
$data1[110,0] = "19:11:59"
$data1[110,1] = "21:02"
$data1[110,2] = "27_EL RETIRO"
$data1[110,3] = 111
$data1[110,4] = "LP1912"
$ws1.Range("A6:E116").Value = $data1

# ---------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = $updatedLabel
$ws2.Range("A3").Value = "Total filas: 15"

$data2 = New-Object 'object[,]' 15,5
$data2[0,0] = "16:46:42"
$data2[0,1] = "16:56"
$data2[0,2] = "215C_EL PATO"
$data2[0,3] = 10
$data2[0,4] = "LP1912"
$data2[1,0] = "16:46:42"
$data2[1,1] = "17:44"
$data2[1,2] = "215B_EL PATO"
$data2[1,3] = 58
$data2[1,4] = "LP1912"
$data2[2,0] = "17:35:09"
$data2[2,1] = "17:45"
$data2[2,2] = "215B_EL PATO"
$data2[2,3] = 10
$data2[2,4] = "LP1912"
$data2[3,0] = "16:46:42"
$data2[3,1] = "17:50"
$data2[3,2] = "215_EL PELIGRO"
$data2[3,3] = 64
$data2[3,4] = "LP1912"
$data2[4,0] = "16:52:37"
$data2[4,1] = "17:51"
$data2[4,2] = "215_EL PELIGRO"
$data2[4,3] = 59
$data2[4,4] = "LP1912"
$data2[5,0] = "17:47:22"
$data2[5,1] = "17:51"
$data2[5,2] = "215B_EL PATO"
$data2[5,3] = 4
$data2[5,4] = "LP1912"
$data2[6,0] = "17:54:43"
$data2[6,1] = "19:20"
$data2[6,2] = "215C_EL PATO"
$data2[6,3] = 86
$data2[6,4] = "LP1912"
$data2[7,0] = "17:35:09"
$data2[7,1] = "19:21"
$data2[7,2] = "215C_EL PATO"
$data2[7,3] = 106
$data2[7,4] = "LP1912"
$data2[8,0] = "17:54:43"
$data2[8,1] = "19:30"
$data2[8,2] = "215_EL PELIGRO"
$data2[8,3] = 96
$data2[8,4] = "LP1912"
$data2[9,0] = "17:35:09"
$data2[9,1] = "19:31"
$data2[9,2] = "215_EL PELIGRO"
$data2[9,3] = 116
$data2[9,4] = "LP1912"
$data2[10,0] = "18:10:41"
$data2[10,1] = "20:00"
$data2[10,2] = "215C_EL PATO"
$data2[10,3] = 110
$data2[10,4] = "LP1912"
$data2[11,0] = "18:31:18"
$data2[11,1] = "20:01"
$data2[11,2] = "215C_EL PATO"
$data2[11,3] = 90
$data2[11,4] = "LP1912"
$data2[12,0] = "19:11:59"
$data2[12,1] = "20:43"
$data2[12,2] = "215B_EL PATO"
$data2[12,3] = 92
$data2[12,4] = "LP1912"
$data2[13,0] = "18:52:04"
$data2[13,1] = "20:44"
$data2[13,2] = "215B_EL PATO"
$data2[13,3] = 112
$data2[13,4] = "LP1912"
$data2[14,0] = "19:11:59"
$data2[14,1] = "21:01"
$data2[14,2] = "215A_EL PATO"
$data2[14,3] = 110
$data2[14,4] = "LP1912"
$ws2.Range("A6:E20").Value = $data2

# ---------------------------------------------------------------
# Sheet 3: 6203-6173 (only the "last updated" header changes)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = $updatedLabel
